# New crime data collected - weekly CompStat update (72nd Precinct)
# Updates the "Volume/Number" and "Week covering" headers, and refreshes
# the crime-complaint statistics table (rows 14-33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text updates (rich-text shared strings): bump the issue number
# and roll the reporting week forward by 7 days.
# ---------------------------------------------------------------------
$ws.Range("A8").Value2 = "Volume 31   Number  38"
$ws.Range("C9").Value2 = "Report Covering the Week  9/16/2024  Through  9/22/2024"

# ---------------------------------------------------------------------
# Row 14 (Murder): the 28-Day columns flip from numeric (1 / -100) to
# the "no data" shared-string markers ("0" / "***.*") used elsewhere in
# the sheet. Copy the format+value from cells that already hold those
# markers (C14 = "0", E14 = "***.*") so the style index matches exactly.
# ---------------------------------------------------------------------
$ws.Range("C14").Copy($ws.Range("G14"))
$ws.Range("E14").Copy($ws.Range("H14"))

# ---------------------------------------------------------------------
# Row 22 (Housing): the Week-to-Date % Chg columns flip the other way,
# from the shared-string markers to real numbers. Borrow number format
# from sibling numeric cells in the same row (C22 / H22), then set the
# real values.
# ---------------------------------------------------------------------
$ws.Range("C22").Copy($ws.Range("D22"))
$ws.Range("H22").Copy($ws.Range("E22"))
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -50

# ---------------------------------------------------------------------
# Row 33 (Hate Crimes): same transformation as row 22, this time for the
# Week-to-Date columns (D33 / E33), borrowing format from the row's own
# 28-day numeric cells (G33 / H33).
# ---------------------------------------------------------------------
$ws.Range("G33").Copy($ws.Range("D33"))
$ws.Range("H33").Copy($ws.Range("E33"))
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = -100

# ---------------------------------------------------------------------
# Remaining plain numeric refreshes across the table (counts and
# percentage changes for rows 14, 16-22, 24-26, 28-30, 33).
# ---------------------------------------------------------------------
$ws.Range("N14").Value = -78.571428571428

$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -75
$ws.Range("F16").Value = 14
$ws.Range("H16").Value = -12.5
$ws.Range("I16").Value = 139
$ws.Range("J16").Value = 113
$ws.Range("K16").Value = 23.008849557522
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -9.740259740259
$ws.Range("N16").Value = -81.342281879194

$ws.Range("F17").Value = 21
$ws.Range("H17").Value = -12.5
$ws.Range("I17").Value = 234
$ws.Range("J17").Value = 216
$ws.Range("K17").Value = 8.333333333333
$ws.Range("L17").Value = 4
$ws.Range("M17").Value = 65.957446808510
$ws.Range("N17").Value = -41.645885286783

$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 25
$ws.Range("I18").Value = 123
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 23
$ws.Range("L18").Value = -35.263157894736
$ws.Range("M18").Value = -34.574468085106
$ws.Range("N18").Value = -84.852216748768

$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -46.153846153846
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = -33.962264150943
$ws.Range("I19").Value = 355
$ws.Range("J19").Value = 414
$ws.Range("K19").Value = -14.251207729468
$ws.Range("L19").Value = -41.89852700491
$ws.Range("M19").Value = 30.03663003663
$ws.Range("N19").Value = 6.606606606606

$ws.Range("G20").Value = 16
$ws.Range("H20").Value = -12.5
$ws.Range("I20").Value = 115
$ws.Range("J20").Value = 128
$ws.Range("K20").Value = -10.15625
$ws.Range("L20").Value = -1.709401709401
$ws.Range("M20").Value = 35.294117647058
$ws.Range("N20").Value = -84.203296703296

$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -18.518518518518
$ws.Range("F21").Value = 99
$ws.Range("G21").Value = 122
$ws.Range("H21").Value = -18.852459016393
$ws.Range("I21").Value = 986
$ws.Range("J21").Value = 988
$ws.Range("K21").Value = -0.202429149797
$ws.Range("L21").Value = -24.095458044649
$ws.Range("M21").Value = 14.385150812065
$ws.Range("N21").Value = -67.735602094240

$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -60
$ws.Range("I22").Value = 14
$ws.Range("J22").Value = 24
$ws.Range("K22").Value = -41.666666666666
$ws.Range("L22").Value = -48.148148148148
$ws.Range("M22").Value = -33.333333333333

$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -45.161290322580
$ws.Range("F24").Value = 68
$ws.Range("G24").Value = 116
$ws.Range("H24").Value = -41.379310344827
$ws.Range("I24").Value = 698
$ws.Range("J24").Value = 914
$ws.Range("K24").Value = -23.632385120350
$ws.Range("L24").Value = -29.919678714859
$ws.Range("M24").Value = 25.314183123877

$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = -69.565217391304
$ws.Range("F25").Value = 19
$ws.Range("G25").Value = 71
$ws.Range("H25").Value = -73.239436619718
$ws.Range("I25").Value = 280
$ws.Range("J25").Value = 510
$ws.Range("K25").Value = -45.098039215686
$ws.Range("L25").Value = -48.998178506375

$ws.Range("C26").Value = 15
$ws.Range("E26").Value = 36.363636363636
$ws.Range("F26").Value = 58
$ws.Range("G26").Value = 40
$ws.Range("H26").Value = 45
$ws.Range("I26").Value = 406
$ws.Range("J26").Value = 394
$ws.Range("K26").Value = 3.045685279187
$ws.Range("L26").Value = 21.556886227544
$ws.Range("M26").Value = -15.062761506276

$ws.Range("D28").Value = 2
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -85.714285714285
$ws.Range("J28").Value = 65
$ws.Range("K28").Value = -35.384615384615
$ws.Range("L28").Value = -17.647058823529

$ws.Range("N29").Value = -90

$ws.Range("N30").Value = -89.473684210526

$ws.Range("G33").Value = 2
$ws.Range("J33").Value = 3
$ws.Range("K33").Value = 66.666666666666
